# Mark a few more strings ("ok") for translation in rows 41-70 of the
# Filelist sheet, and update the sheet view (frozen/split pane + selection)
# to reflect scrolling down to row 61 with B70 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Fill B41:B70 with "ok" (shared string already used throughout column B).
$ws.Range("B41:B70").Value = "ok"

# Recreate the window split/selection state: rows 1-60 frozen/split off,
# top-left cell of the lower pane is A61, and the active cell there is B70.
$win = $excel.ActiveWindow
$win.Split = $true
$win.SplitRow = 60
$win.SplitColumn = 0

$ws.Range("D1:E4").Select()
$ws.Range("B70").Select()
